$d = $word.ActiveDocument

# --- 1. "hdjashdahdasdkk" -> "H" + "djashdahdasdkk" (two runs, same paragraph) ---
# Capitalize the leading "h" first (keeps a single run for now).
$firstChar = $d.Range(0, 1)
$firstChar.Text = "H"

# Split "Hdjashdahdasdkk" into two runs by temporarily inserting a paragraph
# break right after the "H" and then deleting that paragraph mark again.
# Re-joining two paragraphs this way leaves the text in separate <w:r> runs
# instead of Word re-merging them into a single run.
$splitPoint = $d.Range(1, 1)
$splitPoint.InsertParagraphAfter()
$mark = $d.Range(1, 2)
$mark.Delete()

# --- 2. Remove the old bookmark; it will be re-created at the end of the new paragraph ---
$d.Bookmarks("_GoBack").Delete()

# --- 3. Add the new second paragraph with "jJDHjadhjoHDJhdjhdj" ---
$endOfDoc = $d.Content
$endOfDoc.Collapse(0)
$endOfDoc.InsertParagraphAfter()

$newParaStart = $d.Range(16, 16)
$newParaStart.InsertAfter("jJDHjadhjoHDJhdjhdj")

# --- 4. Re-create the "_GoBack" bookmark at the end of the new (second) paragraph ---
# Placing a collapsed bookmark directly on the trailing paragraph-mark position is
# unreliable, so anchor it one character earlier (a safe, interior position),
# then delete+retype the final character so the collapsed bookmark naturally
# drifts forward to sit right after all the paragraph's text.
$bmAnchor = $d.Range(34, 34)
$d.Bookmarks.Add("_GoBack", $bmAnchor)

$tail = $d.Range(34, 35)
$tail.Delete()
$retype = $d.Range(34, 34)
$retype.InsertAfter("j")
